$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 02:51"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 3833134
$ws.Range("C4").Value = 63122
$ws.Range("D4").Value = 1775219
$ws.Range("E4").Value = 1915045
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 806
$ws.Range("H4").Value = 142870

# Row 23: Argentina
$ws.Range("A23").Value = "Argentina"
$ws.Range("B23").Value = 122524
$ws.Range("C23").Value = 3223
$ws.Range("D23").Value = 52607
$ws.Range("E23").Value = 67697
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 42
$ws.Range("H23").Value = 2220

# Row 24: Canada
$ws.Range("A24").Value = "Canada"
$ws.Range("B24").Value = 109999
$ws.Range("C24").Value = 330
$ws.Range("D24").Value = 96914
$ws.Range("E24").Value = 4237
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = 8848

# Row 41: Panama
$ws.Range("A41").Value = "Panama"
$ws.Range("B41").Value = 52261
$ws.Range("C41").Value = 853
$ws.Range("D41").Value = 27494
$ws.Range("E41").Value = 23696
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 33
$ws.Range("H41").Value = 1071

# Row 42: Paises Bajos
$ws.Range("A42").Value = "Paises Bajos"
$ws.Range("B42").Value = 51581
$ws.Range("C42").Value = 127
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 6136

# Row 43: Republica Dominicana
$ws.Range("A43").Value = "Republica Dominicana"
$ws.Range("B43").Value = 51519
$ws.Range("C43").Value = 1406
$ws.Range("D43").Value = 24607
$ws.Range("E43").Value = 25941
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 29
$ws.Range("H43").Value = 971

# Row 48: Guatemala
$ws.Range("A48").Value = "Guatemala"
$ws.Range("B48").Value = 38042
$ws.Range("C48").Value = 4233
$ws.Range("D48").Value = 23365
$ws.Range("E48").Value = 13228
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 6
$ws.Range("H48").Value = 1449

# Row 49: Rumania
$ws.Range("A49").Value = "Rumania"
$ws.Range("B49").Value = 36691
$ws.Range("C49").Value = 889
$ws.Range("D49").Value = 22488
$ws.Range("E49").Value = 12194
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 21
$ws.Range("H49").Value = 2009

# Row 50: Nigeria
$ws.Range("A50").Value = "Nigeria"
$ws.Range("B50").Value = 36107
$ws.Range("C50").Value = 653
$ws.Range("D50").Value = 14938
$ws.Range("E50").Value = 20391
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 778

# Row 51: Barein
$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 36004
$ws.Range("C51").Value = 531
$ws.Range("D51").Value = 31765
$ws.Range("E51").Value = 4115
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 124

# Row 52: Afganistan
$ws.Range("A52").Value = "Afganistan"
$ws.Range("B52").Value = 35301
$ws.Range("C52").Value = 72
$ws.Range("D52").Value = 23273
$ws.Range("E52").Value = 10864
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 17
$ws.Range("H52").Value = 1164

# Row 53: Armenia
$ws.Range("A53").Value = "Armenia"
$ws.Range("B53").Value = 34462
$ws.Range("C53").Value = 461
$ws.Range("D53").Value = 23123
$ws.Range("E53").Value = 10708
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 11
$ws.Range("H53").Value = 631

# Row 70: Chequia
$ws.Range("A70").Value = "Chequia"
$ws.Range("B70").Value = 13855
$ws.Range("C70").Value = 113
$ws.Range("D70").Value = 8733
$ws.Range("E70").Value = 4764
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 358

# Row 128: Suazilandia
$ws.Range("A128").Value = "Suazilandia"
$ws.Range("B128").Value = 1729
$ws.Range("C128").Value = 110
$ws.Range("D128").Value = 788
$ws.Range("E128").Value = 920
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 21

# Row 129: Sierra Leona
$ws.Range("A129").Value = "Sierra Leona"
$ws.Range("B129").Value = 1701
$ws.Range("C129").Value = 13
$ws.Range("D129").Value = 1237
$ws.Range("E129").Value = 399
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 65

# Row 134: Zimbabue
$ws.Range("A134").Value = "Zimbabue"
$ws.Range("B134").Value = 1478
$ws.Range("C134").Value = 58
$ws.Range("D134").Value = 439
$ws.Range("E134").Value = 1014
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 25

# Row 135: Mozambique
$ws.Range("A135").Value = "Mozambique"
$ws.Range("B135").Value = 1435
$ws.Range("C135").Value = 33
$ws.Range("D135").Value = 408
$ws.Range("E135").Value = 1017
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 10

# Row 140: Niger
$ws.Range("A140").Value = "Niger"
$ws.Range("B140").Value = 1104
$ws.Range("C140").Value = 2
$ws.Range("D140").Value = 1014
$ws.Range("E140").Value = 21
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 69

# Row 147: Surinam
$ws.Range("A147").Value = "Surinam"
$ws.Range("B147").Value = 1001
$ws.Range("C147").Value = 58
$ws.Range("D147").Value = 610
$ws.Range("E147").Value = 371
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 20

# Row 148: Republica del Chad
$ws.Range("A148").Value = "Republica del Chad"
$ws.Range("B148").Value = 889
$ws.Range("C148").Value = 2
$ws.Range("D148").Value = 800
$ws.Range("E148").Value = 14
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 75

# Row 162: Vietnam
$ws.Range("A162").Value = "Vietnam"
$ws.Range("B162").Value = 382
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 357
$ws.Range("E162").Value = 25
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 0

# Row 168: Guyana
$ws.Range("A168").Value = "Guyana"
$ws.Range("B168").Value = 327
$ws.Range("C168").Value = 7
$ws.Range("D168").Value = 163
$ws.Range("E168").Value = 145
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 19

# Row 181: Trinidad yTobago
$ws.Range("A181").Value = "Trinidad yTobago"
$ws.Range("B181").Value = 137
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 124
$ws.Range("E181").Value = 5
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 8
